$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.098.97"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.223.63"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "290.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "88.00"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.54%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "30.54"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  -2.03%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.04%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "2.568.69"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "2.209.27"
$ws.Range("E16").Value = "  -1.13%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.731"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "40.045.39"
$ws.Range("E18").Value = "  +0.29%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.52"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +7.61%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -0.84%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.82"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "65.78"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "236.31"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("E24").Value = "  +0.06%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.78%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "22.68"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("E29").Value = "  -0.12%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "155.73"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "31.87"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  +6.30%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "15.79"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.47%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0984"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").Value = "2.116.99"
$ws.Range("E41").Value = "  +8.29%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.83"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.13"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.19%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "10.04"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("E45").Value = "  -1.14%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "17.86"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +9.77%  "
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").Value = "2.434.61"
$ws.Range("E48").Value = "  -0.90%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "89.02"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  -2.64%  "
